# PROS-12983 - CCRU - Wrong KPI calculation
#
# The KPI lists in column Z (Values) and column AA (Scenes to include) used
# comma (",") as the separator between individual list items. Several of the
# downstream KPI values themselves legitimately contain commas (e.g. as part
# of free text), which made a naive comma-split of the cell mis-parse the
# list and produced wrong KPI calculations. The fix re-writes those specific
# cells to use a semicolon (";") as the list separator instead of a comma.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AA ("Scenes to include") - "Panoramic Photo, SS_Panoramic Photo" -> semicolon separated
$aaRows = 4,5,6,7,8,9,10,11,12,13,15,17,18,20,21,23,24,25,26,27,28,29
foreach ($r in $aaRows) {
    $ws.Range("AA$r").Value = "Panoramic Photo; SS_Panoramic Photo"
}

# Column Z (Values) - each distinct comma-separated list rewritten with semicolons
$ws.Range("Z43").Value = "Panoramic photo of Cooler; SS_Panoramic photo of Cooler - Traditional Trade"

$ws.Range("Z46").Value = "SSD 1 door; SSD 1.5 door; Mixed 1 door; Mixed  1.5 door; SS_SSD 1 door; SS_SSD 1.5 door; SS_Mixed 1 door; SS_Mixed  1.5 door"

$z48Value = "SSD 1 door; NCB 1 door; Mixed 1 door; SSD 1.5 door; NSB  1.5 door; Mixed  1.5 door; SSDFL 1 door; SSDFL 1.5 door; SS_SSD 1 door; SS_NCB 1 door; SS_Mixed 1 door; SS_SSD 1.5 door; SS_NSB  1.5 door; SS_Mixed  1.5 door; SS_SSDFL 1 door; SS_SSDFL 1.5 door"
$ws.Range("Z48").Value = $z48Value
$ws.Range("Z49").Value = $z48Value
$ws.Range("Z50").Value = $z48Value

$ws.Range("Z53").Value = "SSD 1 door; NSB 1 door; Mixed 1 door; FC 1 door; SSD 1.5 door; NSB  1.5 door; Mixed  1.5 door; FC door 1.5 door"

$z54Value = "Freezer BIG SSD; SSD Cooler-Side Rack Big; SSD Cooler-Side Rack Small; SSD Mini pallet; SSD Neck-hanging Rack; Freezer Small SSD; Other SSD Display BIG; Other SSD Display Small; Juice Cooler-Side Rack; Juice Freezer; Juice Mini Pallet; Other Juice Display BIG; Pulpy Cooler-Side Rack; Pulpy Other; Other Juice Display Small; Bonaqua Display; Fuzetea Display; Energy Display; Promo SSD; Promo Juice; Promo Water; Promo Energy; Promo Fuzetea"
$ws.Range("Z54").Value = $z54Value
$ws.Range("Z55").Value = $z54Value

# Leave the selection where the edits were made (columns Z:AA), matching the
# reviewer's final on-screen state after making the fix.
$ws.Range("Z1:AA1048576").Select()
